$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 25, shifting all the
# existing data (old rows 25-53) down to rows 28-56.
$ws.Rows("25:27").Insert()

# Fill the 3 newly inserted rows with the new weekly price block
# (same fixed columns as every other data row; only D/L/M/N/O/P/S vary).

# Row 25: Primera
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C25").Value = "Arica y Parinacota"
$ws.Range("D25").Value = 44977
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100101
$ws.Range("H25").Value = "Berries"
$ws.Range("I25").Value = 100112025
$ws.Range("J25").Value = "Frutilla"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 8500
$ws.Range("Q25").Value = "$/bandeja 3 kilos"
$ws.Range("R25").Value = "Región de Arica y Parinacota"
$ws.Range("S25").Value = 2833
$ws.Range("T25").Value = 3

# Row 26: Segunda
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44977
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100112025
$ws.Range("J26").Value = "Frutilla"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 250
$ws.Range("N26").Value = 7000
$ws.Range("O26").Value = 8000
$ws.Range("P26").Value = 7500
$ws.Range("Q26").Value = "$/bandeja 3 kilos"
$ws.Range("R26").Value = "Región de Arica y Parinacota"
$ws.Range("S26").Value = 2500
$ws.Range("T26").Value = 3

# Row 27: Tercera
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44977
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100101
$ws.Range("H27").Value = "Berries"
$ws.Range("I27").Value = 100112025
$ws.Range("J27").Value = "Frutilla"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Tercera"
$ws.Range("M27").Value = 200
$ws.Range("N27").Value = 6000
$ws.Range("O27").Value = 7000
$ws.Range("P27").Value = 6500
$ws.Range("Q27").Value = "$/bandeja 3 kilos"
$ws.Range("R27").Value = "Región de Arica y Parinacota"
$ws.Range("S27").Value = 2167
$ws.Range("T27").Value = 3
